$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19: correct an existing "A" entry in column J to "P"
$ws.Range("J19").Value = "P"

# Fill in previously-blank attendance marks for session 10 (P) and session 11 (Q)
# across every student row. Each row's format is copied from column O of the
# same row first, so the newly-written cells pick up the same "has data" style
# (matching border/alignment) that the rest of the row's filled-in cells already use.
$attendance = [ordered]@{
  7 = @("A", "P")
  8 = @("P", "P")
  9 = @("A", "P")
  10 = @("P", "P")
  11 = @("A", "P")
  12 = @("P", "P")
  13 = @("P", "P")
  14 = @("P", "A")
  15 = @("P", "P")
  16 = @("P", "P")
  17 = @("A", "P")
  18 = @("A", "P")
  19 = @("P", "P")
  20 = @("P", "P")
  21 = @("P", "A")
  22 = @("P", "P")
  23 = @("P", "P")
  24 = @("P", "P")
  25 = @("A", "P")
  26 = @("A", "P")
  27 = @("P", "A")
  28 = @("P", "A")
  29 = @("P", "P")
  30 = @("A", "P")
  31 = @("P", "P")
  32 = @("P", "P")
  33 = @("P", "P")
  34 = @("P", "P")
  35 = @("P", "P")
  36 = @("P", "P")
  37 = @("P", "P")
  38 = @("P", "P")
  39 = @("P", "P")
  40 = @("P", "P")
  41 = @("A", "P")
  42 = @("P", "P")
  43 = @("P", "P")
  44 = @("A", "P")
  45 = @("p", "A")
  46 = @("P", "P")
  47 = @("A", "P")
  48 = @("P", "P")
  49 = @("P", "P")
  50 = @("P", "P")
  51 = @("A", "P")
  52 = @("P", "P")
  53 = @("P", "P")
  54 = @("A", "P")
  55 = @("A", "P")
  56 = @("P", "P")
  57 = @("P", "A")
  58 = @("P", "A")
  59 = @("P", "P")
  60 = @("A", "P")
  61 = @("P", "P")
  62 = @("A", "P")
  63 = @("P", "A")
  64 = @("P", "P")
  65 = @("P", "P")
  66 = @("P", "P")
  67 = @("P", "P")
  68 = @("P", "A")
  69 = @("A", "P")
  70 = @("P", "P")
  71 = @("P", "P")
  72 = @("P", "P")
  73 = @("P", "A")
  74 = @("P", "A")
  75 = @("P", "P")
  76 = @("P", "P")
  77 = @("P", "P")
}

foreach ($r in $attendance.Keys) {
  $vals = $attendance[$r]
  $ws.Range("O$r").Copy()
  $ws.Range("P$r`:Q$r").PasteSpecial(-4122)
  $ws.Range("P$r").Value = $vals[0]
  $ws.Range("Q$r").Value = $vals[1]
}

# Row 84 sits below the last student (row 77) inside the blank template area; only
# its Q cell picked up the "has data" style in the source edit, with no value set.
$ws.Range("O8").Copy()
$ws.Range("Q84").PasteSpecial(-4122)

$excel.CutCopyMode = 0
